$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for B2:O25 (columns B..O), rows 2..25
$rowsData = @(
  @(0.7201429746953067, 0.1347082488026814, 0.0776229439123739, 0.1101984123638573, 0, 0.002506416045459196, 0, 0, 0, 0.3623900200211381, 0.1994137502586426, 0, 0, 5.420126422093233),
  @(0.6818026684744325, 0.1341745273417914, 0.07054334449334476, 0.1095007328195408, 0, 0.002509568969623336, 0, 0, 0, 0.3290667780799197, 0.1926573732315262, 0, 0, 5.418179685006805),
  @(0.6586103868768589, 0.1338413842144561, 0.06623202424644603, 0.1091248149649537, 0, 0.002511608427541074, 0, 0, 0, 0.3087238345941898, 0.18862076084514, 0, 0, 5.419482139380818),
  @(0.6492474953849694, 0.1337042677293709, 0.06448405439245164, 0.1089848260785367, 0, 0.002512465641215367, 0, 0, 0, 0.300463820115823, 0.1870039576222524, 0, 0, 5.420640750839397),
  @(0.6476981315039154, 0.1336814178907382, 0.06419434401105661, 0.1089623785919329, 0, 0.002512609560754462, 0, 0, 0, 0.2990940662048729, 0.1867371896703105, 0, 0, 5.42087105136639),
  @(0.6584837581425518, 0.1338395404994586, 0.06620841438910929, 0.1091228735643064, 0, 0.002511619882217678, 0, 0, 0, 0.3086123156670197, 0.188598842057786, 0, 0, 5.419495222999728),
  @(0.7068510871271769, 0.1345253550723946, 0.07517449903993167, 0.1099469693350166, 0, 0.00250748173002833, 0, 0, 0, 0.3508759216477415, 0.1970609554010991, 0, 0, 5.418936624717503),
  @(0.804453501705126, 0.1358267694415076, 0.0930411140056151, 0.1119790804498422, 0, 0.002500184799597794, 0, 0, 0, 0.4346790871541373, 0.2145423362114229, 0, 0, 5.437676538824746),
  @(0.8778313268151692, 0.1367560605713294, 0.1063448253766524, 0.113725783817948, 0, 0.002495317292733668, 0, 0, 0, 0.4968070698609779, 0.2279282743062225, 0, 0, 5.46356917665679),
  @(0.9115737986624026, 0.1371729189804825, 0.1124363969282172, 0.1145755318736761, 0, 0.002493209019884674, 0, 0, 0, 0.5251910656985217, 0.2341360843832518, 0, 0, 5.477988693911129),
  @(0.924403021388855, 0.1373299194787663, 0.1147488596472925, 0.1149052383757017, 0, 0.002492425828989518, 0, 0, 0, 0.5359566330370171, 0.2365038570485609, 0, 0, 5.483829188265247),
  @(0.9216377273446028, 0.1372961447697065, 0.1142505749243838, 0.1148338778128313, 0, 0.002492593829889732, 0, 0, 0, 0.5336373173514346, 0.2359931587651118, 0, 0, 5.482554420450299),
  @(0.9126282333431845, 0.1371858526872103, 0.1126265299671019, 0.1146024982310152, 0, 0.002493144282609969, 0, 0, 0, 0.5260764132944189, 0.2343305417225423, 0, 0, 5.47846157417564),
  @(0.9071163751936524, 0.1371181840243807, 0.1116325004893213, 0.1144618034593847, 0, 0.002493483423934794, 0, 0, 0, 0.5214473667417963, 0.2333143551305881, 0, 0, 5.476004104139406),
  @(0.8756334333822053, 0.1367286987766789, 0.1059475259850018, 0.1136713602120984, 0, 0.002495457199782536, 0, 0, 0, 0.4949545288145032, 0.2275249610764263, 0, 0, 5.462680000660868),
  @(0.8564121995441951, 0.1364882494574857, 0.1024701439294944, 0.1132005725972185, 0, 0.002496695139190421, 0, 0, 0, 0.4787329555905728, 0.2240036811958106, 0, 0, 5.455182749927133),
  @(0.8453908040917213, 0.1363493963719762, 0.1004737765516097, 0.1129349807312856, 0, 0.002497417149017546, 0, 0, 0, 0.4694142070006251, 0.221989487650248, 0, 0, 5.451119073347854),
  @(0.8416650280248632, 0.1363022884055098, 0.09979848114581102, 0.1128459479611017, 0, 0.002497663325130501, 0, 0, 0, 0.4662610191880674, 0.2213094327910454, 0, 0, 5.449785858110232),
  @(0.8584548015923872, 0.1365139029977271, 0.1028399309483348, 0.1132501513676871, 0, 0.002496562326086647, 0, 0, 0, 0.4804585849273337, 0.2243773733052308, 0, 0, 5.455955118892575),
  @(0.9152731410923991, 0.1372182714135448, 0.1131033960636927, 0.1146702449880443, 0, 0.002492982190335496, 0, 0, 0, 0.5282967708899093, 0.2348184312302806, 0, 0, 5.479653422452373),
  @(0.9527082399392839, 0.1376736307665851, 0.119844492614817, 0.1156445499873726, 0, 0.002490730733556443, 0, 0, 0, 0.5596616849568079, 0.2417413991653632, 0, 0, 5.497357538944755),
  @(0.9327010269330458, 0.1374310560924457, 0.1162435878409127, 0.1151203208862093, 0, 0.002491924316699668, 0, 0, 0, 0.542912608351827, 0.2380374192544252, 0, 0, 5.487705634734198),
  @(0.8575312503809585, 0.1365023069491116, 0.1026727414379422, 0.1132277210122332, 0, 0.002496622338829921, 0, 0, 0, 0.4796784052989267, 0.2242083952147027, 0, 0, 5.455605162764073),
  @(0.7777556220528652, 0.1354793936456069, 0.08817689852158139, 0.1113847912288861, 0, 0.002502071773276481, 0, 0, 0, 0.4119097542123029, 0.2097180236698222, 0, 0, 5.430480055944031)
)

$r = 2
foreach ($row in $rowsData) {
  $c = 2
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value = $val
    $c = $c + 1
  }
  $r = $r + 1
}
